{"js": "// Fix a stray double space in the \"Competencies & Functional Skills:\" line:\n// \"well-rounded,  customer service\" -> \"well-rounded, customer service\"\nconst body = context.document.body;\n\nconst results = body.search(\"well-rounded,  customer service\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"well-rounded, customer service\", \"Replace\");\n  await context.sync();\n} else {\n  // Fallback: search just the doubled space after \"well-rounded,\" in case\n  // surrounding text differs slightly, and collapse it to a single space.\n  const fallback = body.search(\"well-rounded,  customer\", { matchCase: true });\n  fallback.load(\"text\");\n  await context.sync();\n  if (fallback.items.length > 0) {\n    fallback.items[0].insertText(\"well-rounded, customer\", \"Replace\");\n    await context.sync();\n  }\n}\n", "ps1": "# Fix a stray double space in the \"Competencies & Functional Skills:\" line:\n# \"well-rounded,  customer service\" -> \"well-rounded, customer service\"\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"well-rounded,  customer service\"\n$find.Replacement.Text = \"well-rounded, customer service\"\n$found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    # Fallback: collapse just the doubled space, in case the wider phrase\n    # doesn't match exactly.\n    $find2 = $d.Content.Find\n    $find2.ClearFormatting()\n    $find2.Replacement.ClearFormatting()\n    $find2.Text = \"well-rounded,  customer\"\n    $find2.Replacement.Text = \"well-rounded, customer\"\n    $find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n}\n"}
